# Insert a new weekly price-report row for Cilantro ("Agrícola del Norte
# S.A. de Arica") just above the existing row 27. This pushes every
# subsequent record (old rows 27-78) down by one row (to 28-79), matching
# the target diff, and then fills the freshly inserted row 27 with the
# new week's data (market metadata repeated, new date + new price columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 27..78 down to 28..79, leaving a blank row 27 behind.
$ws.Rows("27:27").Insert()

# Populate the new row 27 with the latest weekly observation.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44725
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112040
$ws.Range("G27").Value = "Cilantro"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 2500
$ws.Range("L27").Value = 2800
$ws.Range("M27").Value = 2650
$ws.Range("N27").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 1325
$ws.Range("Q27").Value = 2
$ws.Range("R27").Value = "Hortaliza"
